$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.921.42'
$ws.Range("E2").Value = '  +5.19%  '
$ws.Range("D3").Value = '2.255.40'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.45'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.44%  '
$ws.Range("E7").Value = '  +3.17%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '54.70'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +9.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.23'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.57%  '
$ws.Range("E12").Value = '  +1.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.115'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.64'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").Value = '2.604.68'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.07'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("D17").Value = '2.261.35'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.755'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.11%  '
$ws.Range("D19").Value = '41.795.08'
$ws.Range("E19").Value = '  +4.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.08'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +8.85%  '
$ws.Range("D21").Value = '0.0₃0900'
$ws.Range("E21").Value = '  +1.41%  '
$ws.Range("E22").Value = '  +2.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.98'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.58'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.17%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.89'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.86'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.31'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +12.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.84'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.70'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.44%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +3.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0741'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.77%  '
$ws.Range("E36").Value = '  +1.35%  '
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.115'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.36%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.104'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.41'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +7.13%  '
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.92'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.98%  '
$ws.Range("D43").Value = '2.050.14'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.66'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +8.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0278'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.05'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("E48").Value = '  +4.06%  '
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.60'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.08%  '
